# Bank account unit test plan - completed test data (columns E/F/G, rows 7-22)
# plus a developer-name typo fix in C15 (per author's commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$acctLine = "account_number = 350`nclinet _number = 350`nbalance = 350"

# Row 7 - __init__ / Attributes are set to input values.
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = $acctLine
$ws.Range("G7").Value = "Attribute are created"

# Row 8 - __init__ / Balance attribute set to 0 when non-numeric balance argument.
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = $acctLine
$ws.Range("G8").Value = "balance set to 0"

# Row 9 - __init__ / ValueError when non-numeric account number
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = $acctLine
$ws.Range("G9").Value = "ValueError"

# Row 10 - __init__ / ValueError when non-numeric client number
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = $acctLine
$ws.Range("G10").Value = "ValueError"

# Row 11 - account_number (getter) / returns account number attribute
$ws.Range("E11").Value = $acctLine
$ws.Range("F11").Value = "None"
$ws.Range("G11").Value = "account_number"

# Row 12 - client_number (getter) / returns client number attribute
$ws.Range("E12").Value = $acctLine
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "client_number"

# Row 13 - balance (getter) / returns balance attribute
$ws.Range("E13").Value = $acctLine
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = "balance"

# Row 14 - update_balance / correctly updates balance attribute when positive amount is received.
$ws.Range("E14").Value = $acctLine
$ws.Range("F14").Value = "amount = 100"
$ws.Range("G14").Value = "balance atrribute updated"

# Row 15 - (Method being Tested cell typo'd to developer name) / correctly updates balance attribute when negative amount is received.
$ws.Range("C15").Value = "Sahil Choudhary"
$ws.Range("E15").Value = $acctLine
$ws.Range("F15").Value = "amount = 100"
$ws.Range("G15").Value = "balance atrribute updated"

# Row 16 - update_balance / Balance attribute value remains unchanged when amount is non-numeric
$ws.Range("E16").Value = $acctLine
$ws.Range("F16").Value = 'amount = "Hundred"'
$ws.Range("G16").Value = "balance atrribute updated"

# Row 17 - deposit / BankAccount object's balance is updated correctly when a valid amount is provided.
$ws.Range("E17").Value = $acctLine
$ws.Range("F17").Value = "amount = 100"
$ws.Range("G17").Value = "balance atrribute updated"

# Row 18 - deposit / ValueError when negative amount is provided.
$ws.Range("E18").Value = $acctLine
$ws.Range("F18").Value = "amount = 100"
$ws.Range("G18").Value = "ValueError"

# Row 19 - withdraw / BankAccount object's balance is updated correctly when a valid amount is provided.
$ws.Range("E19").Value = $acctLine
$ws.Range("F19").Value = "amount = 100"
$ws.Range("G19").Value = "balance atrribute updated"

# Row 20 - withdraw / ValueError when negative amount is provided.
$ws.Range("E20").Value = $acctLine
$ws.Range("F20").Value = "amount = -100"
$ws.Range("G20").Value = "ValueError"

# Row 21 - withdraw / ValueError when amount exceeds balance.
$ws.Range("E21").Value = $acctLine
$ws.Range("F21").Value = "amount = 300"
$ws.Range("G21").Value = "ValueError"

# Row 22 - __str__ / returns string in expected format.
$ws.Range("E22").Value = $acctLine
$ws.Range("F22").Value = "None"
$ws.Range("G22").Value = 'Account Number: 350 balance:$350.00'

# Update the saved selection/navigation state to match where the author left off.
$ws.Range("F20").Select()
